$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.364.04"
$ws.Range("E2").Value = "  +7.68%  "
$ws.Range("D3").Value = "3.670.38"
$ws.Range("E3").Value = "  +19.55%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.44%  "
$ws.Range("D7").Value = "3.670.90"
$ws.Range("E7").Value = "  +19.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.537"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.164"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.56"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.500"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +11.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000256"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.50%  "
$ws.Range("D15").Value = "4.290.31"
$ws.Range("E15").Value = "  +19.81%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "71.091.46"
$ws.Range("E16").Value = "  +7.39%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.664.25"
$ws.Range("E17").Value = "  +19.28%  "
$ws.Range("E18").Value = "  +2.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "514.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +21.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.754"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +10.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.90%  "
$ws.Range("E26").Value = "  +10.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.03%  "
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +14.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +16.69%  "
$ws.Range("E32").Value = "  +7.28%  "
$ws.Range("E33").Value = "  +20.23%  "
$ws.Range("E34").Value = "  +6.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.996"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.51%  "
$ws.Range("E37").Value = "  +9.20%  "
$ws.Range("E38").Value = "  +12.25%  "
$ws.Range("B39").Value = "Arweave"
$ws.Range("C39").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "47.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.47%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.13"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.130"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "51.09"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.93"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.77%  "
$ws.Range("D44").Value = "3.155.62"
$ws.Range("E44").Value = "  +13.73%  "
$ws.Range("E45").Value = "  +11.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "407.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +11.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0368"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.95"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +15.65%  "
$ws.Range("E49").Value = "  +16.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "135.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.67%  "
